# PHY410 grade spreadsheet update
# - Raises the max points per homework assignment from 10 to 15.
# - Updates several homework scores to reflect the new 15-point scale.
# - Updates midterm scores/max (now out of 60 instead of 100).
# - Updates the final exam score.
# - Moves the on-screen selection to J13 (scrolled so row 16 is at the top).
# All percentage/average/total formulas recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- HOMEWORK (rows 5-15): points possible 10 -> 15 ----
$ws.Range("C5").Value  = 15
$ws.Range("C6").Value  = 15
$ws.Range("C7").Value  = 15

$ws.Range("B8").Value  = 9
$ws.Range("C8").Value  = 15

$ws.Range("B9").Value  = 14
$ws.Range("C9").Value  = 15

$ws.Range("B10").Value = 15
$ws.Range("C10").Value = 15

$ws.Range("B11").Value = 15
$ws.Range("C11").Value = 15

$ws.Range("B12").Value = 15
$ws.Range("C12").Value = 15

$ws.Range("B13").Value = 15
$ws.Range("C13").Value = 15

$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 15

$ws.Range("B15").Value = 15
$ws.Range("C15").Value = 15

# ---- MIDTERMS (rows 23-25): scores out of 100 -> 60 ----
$ws.Range("B23").Value = 41
$ws.Range("C23").Value = 60

$ws.Range("B24").Value = 44
$ws.Range("C24").Value = 60

$ws.Range("B25").Value = 50
$ws.Range("C25").Value = 60

# ---- FINAL (row 33): score updated ----
$ws.Range("B33").Value = 80

# ---- View state: scroll/select so row 16 is near the top and J13 is active ----
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J13").Select()
